$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on numeric-looking cells so values are stored as text (matching original inline string cells)
$textCells = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","D20","E20","D21","E21","D22","E22","D23","E23","D25","E25","D26","E26","D27","E27","E28","D40","E40","D41","E41","D42","E42","E43","D44","E44","D45","E45","E46","D47","E47","D48","E48","D49","E49","D50","E50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Plain text cell updates (coin names / links)
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("B19").Value = "GateToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("B20").Value = "BTSEToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("B21").Value = "BitpandaEcosystemToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"

# Numeric-looking text cell updates (prices / volumes)
$ws.Range("D2").Value = "244.26"
$ws.Range("E2").Value = "-1.07%"
$ws.Range("D3").Value = "27.19"
$ws.Range("E3").Value = "2.80%"
$ws.Range("D4").Value = "5.099"
$ws.Range("E4").Value = "0.56%"
$ws.Range("D5").Value = "0.05658"
$ws.Range("E5").Value = "0.95%"
$ws.Range("D6").Value = "6.474"
$ws.Range("E6").Value = "-0.45%"
$ws.Range("D7").Value = "0.8202"
$ws.Range("E7").Value = "0.83%"
$ws.Range("D8").Value = "0.8418"
$ws.Range("E8").Value = "0.22%"
$ws.Range("D9").Value = "0.009962"
$ws.Range("E9").Value = "1,569.71%"
$ws.Range("D10").Value = "0.1329"
$ws.Range("E10").Value = "-1.21%"
$ws.Range("D11").Value = "0.06927"
$ws.Range("E11").Value = "-0.52%"
$ws.Range("D12").Value = "0.03158"
$ws.Range("E12").Value = "1.37%"
$ws.Range("D13").Value = "0.02987"
$ws.Range("E13").Value = "5.89%"
$ws.Range("D14").Value = "0.09395"
$ws.Range("E14").Value = "0.09%"
$ws.Range("D15").Value = "0.001512"
$ws.Range("E15").Value = "-0.03%"
$ws.Range("D16").Value = "0.04225"
$ws.Range("E16").Value = "-9.55%"
$ws.Range("D17").Value = "0.006153"
$ws.Range("E17").Value = "0.06%"
$ws.Range("D18").Value = "3.514"
$ws.Range("E18").Value = "-1.14%"
$ws.Range("D19").Value = "3.003"
$ws.Range("E19").Value = "-1.44%"
$ws.Range("D20").Value = "2.307"
$ws.Range("E20").Value = "8.90%"
$ws.Range("D21").Value = "0.3113"
$ws.Range("E21").Value = "-2.14%"
$ws.Range("D22").Value = "0.1254"
$ws.Range("E22").Value = "-3.57%"
$ws.Range("D23").Value = "3.552"
$ws.Range("E23").Value = "-4.96%"
$ws.Range("D25").Value = "0.001227"
$ws.Range("E25").Value = "-1.70%"
$ws.Range("D26").Value = "0.004458"
$ws.Range("E26").Value = "-3.39%"
$ws.Range("D27").Value = "0.00009802"
$ws.Range("E27").Value = "2.17%"
$ws.Range("E28").Value = "-0.02%"
$ws.Range("D40").Value = "0.03656"
$ws.Range("E40").Value = "-0.20%"
$ws.Range("D41").Value = "0.006051"
$ws.Range("E41").Value = "77.20%"
$ws.Range("D42").Value = "0.1052"
$ws.Range("E42").Value = "-22.61%"
$ws.Range("E43").Value = "-11.32%"
$ws.Range("D44").Value = "0.008284"
$ws.Range("E44").Value = "-1.96%"
$ws.Range("D45").Value = "0.00005318"
$ws.Range("E45").Value = "0.49%"
$ws.Range("E46").Value = "0.08%"
$ws.Range("D47").Value = "0.1282"
$ws.Range("E47").Value = "-19.76%"
$ws.Range("D48").Value = "0.002819"
$ws.Range("E48").Value = "36.75%"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("E49").Value = "0.08%"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("E50").Value = "0.08%"
